# Insert a new data row at row 467 (pushing existing rows 467..506 down to 468..507)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 467, shifting rows 467-506 -> 468-507
$ws.Rows.Item(467).Insert()

# Fill in the new row 467 with the new record
$ws.Range("A467").Value = 3
$ws.Range("B467").Value = "Femacal de La Calera"
$ws.Range("C467").Value = "Coquimbo"
$ws.Range("D467").Value = 44769
$ws.Range("E467").Value = 5
$ws.Range("F467").Value = 100112021
$ws.Range("G467").Value = "Ají"
$ws.Range("H467").Value = "Inferno"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 76
$ws.Range("K467").Value = 14000
$ws.Range("L467").Value = 15000
$ws.Range("M467").Value = 14500
$ws.Range("N467").Value = "$/caja 15 kilos"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 967
$ws.Range("Q467").Value = 15
$ws.Range("R467").Value = "Hortaliza"
